# v0.5: Added Ingestion pipeline pseudocode and its short explanation.
#
# Row 2 of the sheet corresponds to the "Ingestion pipeline" story /
# "Write ingestion pseudocode and architecture plan" task. Its "Actual"
# (column E) and "Comment" (column F) cells were previously blank; fill
# them in with the actual time spent and a short explanation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "30 minutes"
$ws.Range("F2").Value = "Pseudo code for ingestion pipeline and explanation is give in the ingestion_pipeline markdown file."
